$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 6.023199999999997
$ws.Range("A3").Value = -21.45310000000003
$ws.Range("B5").Value = 5.020500000000002
$ws.Range("A14").Value = -20.47839999999998
$ws.Range("A16").Value = -21.40740000000002
$ws.Range("B16").Value = 5.4681
$ws.Range("A21").Value = -21.02879999999999
$ws.Range("A23").Value = -21.39700000000003
$ws.Range("A25").Value = -22.37090000000004
